$d = $word.ActiveDocument

$find = "Perioadele campaniei din Constelația Bootes 2022: 14-23 mai, 13-22 iunie, 12-21 iulie"
$replace = "Perioadele campaniei din 2022 pentru Constelația Bootes: 14-23 mai, 13-22 iunie, 12-21 iulie"

$range = $d.Content
$range.Find.Execute($find, $true, $false, $false, $false, $false, `
                     $true, 1, $false, $replace, 2)
